# Updated cryptos list on Tue Dec 19 07:35:30 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.007.71"
$ws.Range("E2").Value = "  +4.00%  "
$ws.Range("D3").Value = "2.246.65"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.36"
$ws.Range("E5").Value = "  +3.02%  "
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.07"
$ws.Range("E7").Value = "  +8.34%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.617"
$ws.Range("E9").Value = "  +6.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.96"
$ws.Range("E10").Value = "  +2.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0936"
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.96"
$ws.Range("E12").Value = "  +2.93%  "
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").Value = "2.586.36"
$ws.Range("E14").Value = "  +3.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.62"
$ws.Range("E15").Value = "  +4.31%  "
$ws.Range("D16").Value = "2.246.67"
$ws.Range("E16").Value = "  +2.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.807"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("D18").Value = "42.916.40"
$ws.Range("E18").Value = "  +4.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.17"
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.99"
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.03"
$ws.Range("E22").Value = "  +3.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.91"
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.17"
$ws.Range("E24").Value = "  +12.33%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.87"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.43"
$ws.Range("E27").Value = "  -2.67%  "
$ws.Range("E28").Value = "  +2.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.74"
$ws.Range("E29").Value = "  +25.13%  "
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.51"
$ws.Range("E31").Value = "  +3.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.33"
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0796"
$ws.Range("E33").Value = "  +3.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.34"
$ws.Range("E34").Value = "  +3.69%  "
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("E36").Value = "  +8.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.34"
$ws.Range("E37").Value = "  +5.00%  "
$ws.Range("E38").Value = "  +17.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.96"
$ws.Range("E39").Value = "  +10.83%  "
$ws.Range("E40").Value = "  +2.24%  "
$ws.Range("E41").Value = "  +2.05%  "
$ws.Range("E42").Value = "  +6.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "60.06"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "105.09"
$ws.Range("E44").Value = "  +6.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.70"
$ws.Range("E45").Value = "  +5.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0994"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.460"
$ws.Range("E47").Value = "  +22.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.39"
$ws.Range("E48").Value = "  +7.98%  "
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("D51").Value = "2.460.66"
$ws.Range("E51").Value = "  +3.17%  "
